$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# Determine the next free row right after the last populated one (row 30 -> 31)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Add a new localization entry: Key / English / German / Polish
$ws.Cells.Item($newRow, 1).Value = "LevelCompleteTitle"
$ws.Cells.Item($newRow, 2).Value = "Great!"
$ws.Cells.Item($newRow, 3).Value = "Super!"
$ws.Cells.Item($newRow, 4).Value = "Super!"

$ws.Cells.Item($newRow, 4).Select()
